# Apply cryptos list update (prices + 1h volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.814.84"
$ws.Range("E2").Value = "  -5.77%  "
$ws.Range("D3").Value = "3.287.74"
$ws.Range("E3").Value = "  -6.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.70"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.95"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Value = "3.280.89"
$ws.Range("E9").Value = "  -6.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  -9.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.584"
$ws.Range("E11").Value = "  -6.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.35"
$ws.Range("E12").Value = "  -8.23%  "
$ws.Range("E13").Value = "  -7.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "643.68"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").Value = "3.804.81"
$ws.Range("E16").Value = "  -6.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.07"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "65.776.05"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "3.283.87"
$ws.Range("E20").Value = "  -6.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.35"
$ws.Range("E21").Value = "  -8.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.903"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.39"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "108.09"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.92"
$ws.Range("E25").Value = "  -8.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").Value = "  -7.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -7.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.67"
$ws.Range("E29").Value = "  -7.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.31"
$ws.Range("E30").Value = "  -7.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.94"
$ws.Range("E31").Value = "  -7.47%  "
$ws.Range("E32").Value = "  -6.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.05"
$ws.Range("E33").Value = "  -5.35%  "
$ws.Range("D35").Value = "3.786.02"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.54"
$ws.Range("E36").Value = "  -6.50%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "518.88"
$ws.Range("E38").Value = "  -6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  -6.59%  "
$ws.Range("D40").Value = "0.0₃0734"
$ws.Range("E40").Value = "  -7.23%  "
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -6.22%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "32.97"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  -18.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  -10.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.22"
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  -8.71%  "
$ws.Range("E51").Value = "  +1.88%  "
